$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 10 (only across columns A:D so we don't disturb the
# rest of the sheet / create stray full-row formatting), shifting
# rows 10-20 down to 11-21.
$ws.Range("A10:D10").Insert(-4121)

# New row 10 should look like row 3 (label/description/value layout):
# columns A and C left-aligned, columns B and D centered, all with the
# same thin-bordered "table cell" style used throughout the sheet.
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("A10").Value2 = "user"
$ws.Range("C10").Value2 = "Go back to user menu"

# Selection moved from B8 to C8
$ws.Range("C8").Select()
